$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G5:H6").NumberFormat = "0.00"
$ws.Range("F7").Copy()
$ws.Range("G5:H6").PasteSpecial(-4122)
